# Enhance shopkeeper management: add brand field to shopkeeper queries;
# implement product fetching by brand and add shopkeeper payment details retrieval
#
# Concretely this applies the following changes to the "Orders" sheet:
#  - Rows 34-43 (column G / Order_Info) get the text "Order Description"
#  - Row 47 gets numeric Shopkeeper_ID (B47=2) and Salesman_ID (C47=1)
#  - A brand-new row 49 is appended with order data

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Order Description" text for column G in rows 34 through 43
for ($r = 34; $r -le 43; $r++) {
    $ws.Cells.Item($r, 7).Value = "Order Description"
}

# Row 47 was missing Shopkeeper_ID / Salesman_ID values - set them now
$ws.Range("B47").Value = 2
$ws.Range("C47").Value = 1

# Append a brand-new order row (row 49)
$ws.Range("A49").Value = 48
$ws.Range("B49").Value = 7
$ws.Range("C49").Value = 2
$ws.Range("D49").Value = "2025-03-16 19:02:52"
$ws.Range("E49").Value = 2830.326799810271
$ws.Range("F49").Value = 0
$ws.Range("G49").Value = "Order 1202"
